$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "bneq" Branch control bit bug (row 11)
$ws.Range("G11").Value = 1

# Update the derived ControlBin and ControlHex values for row 11
$ws.Range("O11").Value = "10010001100000"
$ws.Range("P11").Value = 2460

# Update the selected cell to reflect post-edit selection
$ws.Range("O18").Select()
